# Add a new "iowa_k" parameter sheet to the workbook.
#
# The new sheet is a duplicate of "high_k" (same layout/values) but with
# the permeability ("k") average changed to 3, representing a new well
# diameter / k sweep case. The new sheet is appended as the last tab and
# becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Duplicate the "high_k" sheet and place the copy after the last sheet.
$source = $wb.Worksheets.Item("high_k")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

# The copy is now the last worksheet - rename it.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "iowa_k"

# Update the permeability (k) average value for this sweep case.
$newSheet.Range("C5").Value = 3

# Match the author's last selection on the new sheet.
$newSheet.Range("C6").Select() | Out-Null
